$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Z column (Ring Used) values for rows 2 and 3 to text references "A" and "B"
$ws.Range("Z2").Value = "A"
$ws.Range("Z3").Value = "B"

# Update the view state: scroll to O1, select Z8
$ws.Application.ActiveWindow.ScrollColumn = 15
$ws.Range("Z8").Select()
